# Update cryptos list with refreshed prices / 1h volume percentages,
# and apply the two rank swaps (InternetComputer<->Maker at rows 33/34,
# BabyDogeCoin<->Algorand at rows 49/50).
# Leading '' on purely numeric-looking price strings forces Excel to keep
# them as text (preserving exact formatting, e.g. trailing zeros like
# "210.96" or "65.30", instead of being reinterpreted as a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.867.03'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.624.03'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''210.96'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = '''23.46'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.854.37'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.617.36'
$ws.Range('E13').Value = '  -1.34%  '
$ws.Range('E14').Value = '  -1.93%  '
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = '''65.30'
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').Value = '27.853.79'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = '''229.45'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').Value = '''7.64'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('E23').Value = '  -5.50%  '
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('D25').Value = '''154.73'
$ws.Range('E25').Value = '  +2.04%  '
$ws.Range('D26').Value = '''6.89'
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = '''15.51'
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''3.07'
$ws.Range('E33').Value = '  -1.19%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.398.42'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('E36').Value = '  +9.16%  '
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').Value = '''0.555'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '''65.88'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = '1.764.94'
$ws.Range('E47').Value = '  -0.96%  '
$ws.Range('D48').Value = '''87.89'
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.102'
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('E51').Value = '  -0.57%  '
